$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A 20,000 amount that used to be folded into columns F..K of the "tot" row
# (row 17) is pulled back out into its own cell (F21), so each of those six
# monthly totals is reduced by 20,000.
$ws.Range("L17:M17").Formula = "=+L15+L16"
$ws.Range("F17").Formula = "=(+F15+F16)-20000"
$ws.Range("G17").Formula = "=(+G15+G16)-20000"
$ws.Range("H17").Formula = "=(+H15+H16)-20000"
$ws.Range("I17").Formula = "=(+I15+I16)-20000"
$ws.Range("J17").Formula = "=(+J15+J16)-20000"
$ws.Range("K17").Formula = "=(+K15+K16)-20000"

# The 20,000 that was carved out now lives by itself a few rows down.
$ws.Range("F21").Value = 20000

# Reflect the author's new active selection.
$ws.Range("I7").Select()
